# Mise à jour de l'application
# Adds 13 new training-log rows (192-204) for the 2025-08-26 session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# --- new data rows -------------------------------------------------------
# Columns: Date(A) Nom(B) Volume(C) Intensite(D) Fatigue(E) Douleur(F)
#          Localisation douleur(G) Plaisir(H) Charge(I, = C*D)
$rows = @(
    @{ R=192; Nom="Romain Thunet";    C=48; D=5; E=8; F=4; G="Orteil";            H=0 },
    @{ R=193; Nom="Emmanuel Valey";   C=63; D=1; E=1; F=5; G="Adducteur";         H=6 },
    @{ R=194; Nom="Amir Etien";       C=48; D=2; E=2; F=1; G="Courbature";        H=4 },
    @{ R=195; Nom="Yoan Zouma";       C=48; D=5; E=8; F=5; G="";                  H=4 },
    @{ R=196; Nom="Omar Benyounes";   C=68; D=5; E=5; F=0; G="";                  H=6 },
    @{ R=197; Nom="Naim Ighbane";     C=68; D=7; E=3; F=2; G=("Mollet" + $nbsp);  H=3 },
    @{ R=198; Nom="Yanis Berrached";  C=48; D=5; E=9; F=0; G="";                  H=5 },
    @{ R=199; Nom="Amine Taiar";      C=68; D=2; E=4; F=7; G=("Genou" + $nbsp);   H=8 },
    @{ R=200; Nom="Ilan Ihaddadene";  C=48; D=6; E=6; F=0; G="";                  H=6 },
    @{ R=201; Nom="Hedi Nasri";       C=68; D=6; E=5; F=3; G="Ischio";            H=6 },
    @{ R=202; Nom="Wael Fareh";       C=68; D=7; E=6; F=0; G="";                  H=8 },
    @{ R=203; Nom="Sofiane Belle";    C=48; D=3; E=5; F=0; G="";                  H=2 },
    @{ R=204; Nom="Naim Dhib";        C=48; D=3; E=6; F=0; G="";                  H=3 }
)

$newDate = 45895

foreach ($row in $rows) {
    $r = $row.R

    # Copy the formatting (styles) of the last existing data row (191) down
    # onto the new row, so date/text/number formatting matches the rest of
    # the table.
    $ws.Range("A191:I191").Copy() | Out-Null
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Range("A$r").Value = $newDate
    $ws.Range("B$r").Value = $row.Nom
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    if ($row.G -ne "") {
        # Cells that carry a "Localisation douleur" value use the other
        # font style (same as e.g. G4) rather than the empty-cell style
        # that got pasted from row 191's (empty) G cell above.
        $ws.Range("G4").Copy() | Out-Null
        $ws.Range("G$r").PasteSpecial(-4122) | Out-Null
        $ws.Range("G$r").Value = $row.G
    }
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Formula = "=C$r*D$r"
}

$excel.CutCopyMode = $false

# --- restore the view scroll position / selection recorded in the diff ---
$ws.Range("L197").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 176
$excel.ActiveWindow.ScrollColumn = 1
